$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Cell, $Value)
    $Cell.NumberFormat = "@"
    $Cell.Value = $Value
    $Cell.Style = "Normal"
}

# Row 2
Set-TextValue $ws.Range('D2') '331.53'
Set-TextValue $ws.Range('E2') '0.01%'
Set-TextValue $ws.Range('G2') '17'

# Row 3
Set-TextValue $ws.Range('D3') '41.52'
Set-TextValue $ws.Range('E3') '-0.69%'
Set-TextValue $ws.Range('G3') '17'

# Row 4
Set-TextValue $ws.Range('D4') '5.645'
Set-TextValue $ws.Range('E4') '-1.17%'
Set-TextValue $ws.Range('G4') '17'

# Row 5
Set-TextValue $ws.Range('D5') '0.08341'
Set-TextValue $ws.Range('E5') '2.80%'
Set-TextValue $ws.Range('G5') '17'

# Row 6
Set-TextValue $ws.Range('D6') '8.785'
Set-TextValue $ws.Range('E6') '0.53%'
Set-TextValue $ws.Range('G6') '17'

# Row 7
Set-TextValue $ws.Range('D7') '1.990'
Set-TextValue $ws.Range('E7') '-3.55%'
Set-TextValue $ws.Range('G7') '17'

# Row 8
Set-TextValue $ws.Range('D8') '4.491'
Set-TextValue $ws.Range('E8') '-0.78%'
Set-TextValue $ws.Range('G8') '17'

# Row 9
Set-TextValue $ws.Range('D9') '2.922'
Set-TextValue $ws.Range('E9') '-1.70%'
Set-TextValue $ws.Range('G9') '17'

# Row 10
Set-TextValue $ws.Range('D10') '0.9255'
Set-TextValue $ws.Range('E10') '0.02%'
Set-TextValue $ws.Range('G10') '17'

# Row 11
Set-TextValue $ws.Range('D11') '0.1290'
Set-TextValue $ws.Range('E11') '2.26%'
Set-TextValue $ws.Range('G11') '17'

# Row 12
Set-TextValue $ws.Range('D12') '0.1973'
Set-TextValue $ws.Range('E12') '0.73%'
Set-TextValue $ws.Range('G12') '17'

# Row 13
Set-TextValue $ws.Range('D13') '0.09533'
Set-TextValue $ws.Range('E13') '3.76%'
Set-TextValue $ws.Range('G13') '17'

# Row 14
Set-TextValue $ws.Range('D14') '0.03910'
Set-TextValue $ws.Range('E14') '5.93%'
Set-TextValue $ws.Range('G14') '17'

# Row 15
Set-TextValue $ws.Range('E15') '0.53%'
Set-TextValue $ws.Range('G15') '17'

# Row 16
Set-TextValue $ws.Range('D16') '0.001309'
Set-TextValue $ws.Range('E16') '-0.02%'
Set-TextValue $ws.Range('G16') '17'

# Row 17
Set-TextValue $ws.Range('D17') '0.006107'
Set-TextValue $ws.Range('E17') '-0.05%'
Set-TextValue $ws.Range('G17') '17'

# Row 18
Set-TextValue $ws.Range('B18') 'HotbitToken'
Set-TextValue $ws.Range('C18') 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'
Set-TextValue $ws.Range('D18') '0.004382'
Set-TextValue $ws.Range('E18') '-2.41%'
Set-TextValue $ws.Range('G18') '17'

# Row 19
Set-TextValue $ws.Range('B19') 'LEO'
Set-TextValue $ws.Range('C19') 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
Set-TextValue $ws.Range('D19') '3.441'
Set-TextValue $ws.Range('E19') '1.82%'
Set-TextValue $ws.Range('G19') '17'

# Row 20
Set-TextValue $ws.Range('B20') 'BitpandaEcosystemToken'
Set-TextValue $ws.Range('C20') 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
Set-TextValue $ws.Range('D20') '0.3538'
Set-TextValue $ws.Range('E20') '0.29%'
Set-TextValue $ws.Range('G20') '17'

# Row 21
Set-TextValue $ws.Range('B21') 'MCDex'
Set-TextValue $ws.Range('C21') 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
Set-TextValue $ws.Range('D21') '8.523'
Set-TextValue $ws.Range('E21') '-3.43%'
Set-TextValue $ws.Range('G21') '17'

# Row 22
Set-TextValue $ws.Range('B22') 'ProBitToken'
Set-TextValue $ws.Range('C22') 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'
Set-TextValue $ws.Range('D22') '0.1372'
Set-TextValue $ws.Range('E22') '-3.13%'
Set-TextValue $ws.Range('G22') '17'

# Row 23
Set-TextValue $ws.Range('B23') 'ZBToken'
Set-TextValue $ws.Range('C23') 'https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb'
Set-TextValue $ws.Range('D23') '0.2443'
Set-TextValue $ws.Range('E23') '-6.43%'
Set-TextValue $ws.Range('G23') '17'

# Row 24
Set-TextValue $ws.Range('B24') 'CoinExToken'
Set-TextValue $ws.Range('C24') 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
Set-TextValue $ws.Range('D24') '0.04407'
Set-TextValue $ws.Range('E24') '-0.61%'
Set-TextValue $ws.Range('G24') '17'

# Row 25
Set-TextValue $ws.Range('B25') 'BitKan'
Set-TextValue $ws.Range('C25') 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'
Set-TextValue $ws.Range('D25') '0.001274'
Set-TextValue $ws.Range('E25') '1.31%'
Set-TextValue $ws.Range('G25') '17'

# Row 26
Set-TextValue $ws.Range('D26') '0.0001199'
Set-TextValue $ws.Range('E26') '-3.43%'
Set-TextValue $ws.Range('G26') '17'

# Row 39
Set-TextValue $ws.Range('D39') '0.02819'
Set-TextValue $ws.Range('E39') '0.73%'
Set-TextValue $ws.Range('G39') '17'

# Row 40
Set-TextValue $ws.Range('D40') '0.05514'
Set-TextValue $ws.Range('E40') '0.04%'
Set-TextValue $ws.Range('G40') '17'

# Row 41
Set-TextValue $ws.Range('D41') '0.007956'
Set-TextValue $ws.Range('E41') '2.72%'
Set-TextValue $ws.Range('G41') '17'

# Row 42
Set-TextValue $ws.Range('D42') '0.1438'
Set-TextValue $ws.Range('E42') '1.21%'
Set-TextValue $ws.Range('G42') '17'

# Row 43
Set-TextValue $ws.Range('D43') '0.009326'
Set-TextValue $ws.Range('E43') '-6.04%'
Set-TextValue $ws.Range('G43') '17'

# Row 44
Set-TextValue $ws.Range('D44') '0.002139'
Set-TextValue $ws.Range('E44') '-3.40%'
Set-TextValue $ws.Range('G44') '17'

# Row 45
Set-TextValue $ws.Range('D45') '0.01105'
Set-TextValue $ws.Range('E45') '4.00%'
Set-TextValue $ws.Range('G45') '17'

# Row 46
Set-TextValue $ws.Range('D46') '0.00007102'
Set-TextValue $ws.Range('E46') '4.24%'
Set-TextValue $ws.Range('G46') '17'

# Row 47
Set-TextValue $ws.Range('D47') '0.00000000750'
Set-TextValue $ws.Range('E47') '-0.21%'
Set-TextValue $ws.Range('G47') '17'

# Row 48
Set-TextValue $ws.Range('D48') '0.003236'
Set-TextValue $ws.Range('E48') '7.99%'
Set-TextValue $ws.Range('G48') '17'

# Row 49
Set-TextValue $ws.Range('D49') '0.002279'
Set-TextValue $ws.Range('E49') '-0.15%'
Set-TextValue $ws.Range('G49') '17'

# Row 50
Set-TextValue $ws.Range('D50') '0.00002099'
Set-TextValue $ws.Range('E50') '-0.21%'
Set-TextValue $ws.Range('G50') '17'

# Row 51
Set-TextValue $ws.Range('D51') '0.0001999'
Set-TextValue $ws.Range('E51') '-0.21%'
Set-TextValue $ws.Range('G51') '17'

# Rows with only the Hora (G) column change
Set-TextValue $ws.Range('G27') '17'
Set-TextValue $ws.Range('G28') '17'
Set-TextValue $ws.Range('G29') '17'
Set-TextValue $ws.Range('G30') '17'
Set-TextValue $ws.Range('G31') '17'
Set-TextValue $ws.Range('G32') '17'
Set-TextValue $ws.Range('G33') '17'
Set-TextValue $ws.Range('G34') '17'
Set-TextValue $ws.Range('G35') '17'
Set-TextValue $ws.Range('G36') '17'
Set-TextValue $ws.Range('G37') '17'
Set-TextValue $ws.Range('G38') '17'
